$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Bug fix #1: the onboard-LED enable/disable values (G2:G9) must no longer
#    be hard coded in the sheet - clear their contents but keep formatting.
# ---------------------------------------------------------------------------
$ws.Range("G2:G9").ClearContents()

# ---------------------------------------------------------------------------
# 2. Bug fix #2: document the new "load_radar_parameter" flash-restore entry.
#    Insert a brand new row above row 16 (pushes the existing rows 16-22 down
#    to 17-23) and fill it in with the new parameter description.
# ---------------------------------------------------------------------------
$ws.Rows(16).Insert()

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "load_radar_parameter"
$ws.Range("C16").Value = "u32"
$ws.Range("D16").Value = "加载覆盖范围"
$ws.Range("E16").Value = "UPSSA0"
$ws.Range("F16").Value = "0x38"
$ws.Range("G16").Value = 1

# Re-number the rows that got shifted down by the insert above (their "序号"
# column held literal numbers that need to move on by one).
$ws.Range("A17").Value = 16
$ws.Range("A18").Value = 17
$ws.Range("A19").Value = 18

# ---------------------------------------------------------------------------
# Leave the cursor where the author left it after making the edit.
# ---------------------------------------------------------------------------
$ws.Range("C24").Select()
